{"js": "const pairs = [\n  [\"44\u00d780=3520\", \"39\u00d753=2067\"],\n  [\"62\u00d745=2790\", \"41\u00d779=3239\"],\n  [\"95\u00d788=8360\", \"28\u00d748=1344\"],\n  [\"83\u00d787=7221\", \"37\u00d741=1517\"],\n  [\"31\u00d754=1674\", \"77\u00d711=847\"],\n  [\"60\u00d769=4140\", \"77\u00d781=6237\"],\n  [\"64\u00d728=1792\", \"24\u00d752=1248\"],\n  [\"11\u00d787=957\", \"34\u00d775=2550\"],\n  [\"42\u00d791=3822\", \"56\u00d792=5152\"],\n  [\"11\u00d737=407\", \"38\u00d756=2128\"],\n  [\"50\u00d718=900\", \"48\u00d757=2736\"],\n  [\"62\u00d767=4154\", \"31\u00d731=961\"],\n  [\"79\u00d718=1422\", \"51\u00d747=2397\"],\n  [\"74\u00d726=1924\", \"70\u00d769=4830\"],\n  [\"61\u00d775=4575\", \"42\u00d726=1092\"],\n  [\"77\u00d741=3157\", \"49\u00d738=1862\"],\n  [\"86\u00d759=5074\", \"91\u00d760=5460\"],\n  [\"77\u00d762=4774\", \"44\u00d722=968\"],\n  [\"42\u00d775=3150\", \"30\u00d778=2340\"],\n  [\"93\u00d721=1953\", \"11\u00d746=506\"],\n  [\"42\u00d777=3234\", \"71\u00d747=3337\"],\n  [\"74\u00d731=2294\", \"74\u00d765=4810\"],\n  [\"17\u00d766=1122\", \"27\u00d725=675\"],\n  [\"15\u00d797=1455\", \"78\u00d723=1794\"],\n  [\"93\u00d730=2790\", \"85\u00d730=2550\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('44\u00d780=3520', '39\u00d753=2067'),\n    @('62\u00d745=2790', '41\u00d779=3239'),\n    @('95\u00d788=8360', '28\u00d748=1344'),\n    @('83\u00d787=7221', '37\u00d741=1517'),\n    @('31\u00d754=1674', '77\u00d711=847'),\n    @('60\u00d769=4140', '77\u00d781=6237'),\n    @('64\u00d728=1792', '24\u00d752=1248'),\n    @('11\u00d787=957', '34\u00d775=2550'),\n    @('42\u00d791=3822', '56\u00d792=5152'),\n    @('11\u00d737=407', '38\u00d756=2128'),\n    @('50\u00d718=900', '48\u00d757=2736'),\n    @('62\u00d767=4154', '31\u00d731=961'),\n    @('79\u00d718=1422', '51\u00d747=2397'),\n    @('74\u00d726=1924', '70\u00d769=4830'),\n    @('61\u00d775=4575', '42\u00d726=1092'),\n    @('77\u00d741=3157', '49\u00d738=1862'),\n    @('86\u00d759=5074', '91\u00d760=5460'),\n    @('77\u00d762=4774', '44\u00d722=968'),\n    @('42\u00d775=3150', '30\u00d778=2340'),\n    @('93\u00d721=1953', '11\u00d746=506'),\n    @('42\u00d777=3234', '71\u00d747=3337'),\n    @('74\u00d731=2294', '74\u00d765=4810'),\n    @('17\u00d766=1122', '27\u00d725=675'),\n    @('15\u00d797=1455', '78\u00d723=1794'),\n    @('93\u00d730=2790', '85\u00d730=2550'),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
